$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14:D14").Copy()
$ws.Range("B15:D15").PasteSpecial(-4122)
$ws.Range("B14:D14").Copy()
$ws.Range("B16:D16").PasteSpecial(-4122)

$ws.Cells.Item(15, 2).Value2 = 45333
$ws.Cells.Item(15, 3).Value2 = 7
$ws.Cells.Item(15, 4).Value2 = "Tein logiikan muuntimelle ja parantelin muunnin sivun ulkonäköä."

$ws.Cells.Item(16, 2).Value2 = 45334
$ws.Cells.Item(16, 3).Value2 = 10
$ws.Cells.Item(16, 4).Value2 = "Kehitin verko töykaluihin pinger työkalun ja myöskin kehitin lunttilappu sivun rakenteen ja toiminallisuuden."

$ws.Range("B15:D15").Copy()
$ws.Range("B17:D17").PasteSpecial(-4122)

$ws.Cells.Item(17, 2).Value2 = "Yht"
$ws.Cells.Item(18,1).Formula = "=SUM(C6:C14)"
$ws.Cells.Item(19,1).Formula = "=SUM(C15:C16)"
$ws.Cells.Item(20,1).Formula = "=SUM(C6:C16)"
$ws.Cells.Item(21,1).Formula = "=C15+C16"
"C18 (6:14): " + $ws.Cells.Item(18,1).Value2
"C19 (15:16): " + $ws.Cells.Item(19,1).Value2
"C20 (6:16): " + $ws.Cells.Item(20,1).Value2
"C21 (15+16): " + $ws.Cells.Item(21,1).Value2

$ws.Cells.Item(17, 3).Formula = "=SUM(C6:C16)"
"C17: " + $ws.Cells.Item(17,3).Value2
